$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 429, pushing the old
# rows 429-447 down to become rows 431-449 (all their data, including
# dates/prices/units/origin, moves down untouched).
$ws.Rows.Item(429).Insert()
$ws.Rows.Item(429).Insert()

# Populate the first new row (429) with the latest weekly price report.
$ws.Cells.Item(429, 1).Value2  = 4
$ws.Cells.Item(429, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(429, 3).Value2  = "Los Lagos"
$ws.Cells.Item(429, 4).Value2  = 45041
$ws.Cells.Item(429, 5).Value2  = 10
$ws.Cells.Item(429, 6).Value2  = 100112040
$ws.Cells.Item(429, 7).Value2  = "Cilantro"
$ws.Cells.Item(429, 8).Value2  = "Sin especificar"
$ws.Cells.Item(429, 9).Value2  = "Primera"
$ws.Cells.Item(429, 10).Value2 = 80
$ws.Cells.Item(429, 11).Value2 = 13000
$ws.Cells.Item(429, 12).Value2 = 13000
$ws.Cells.Item(429, 13).Value2 = 13000
$ws.Cells.Item(429, 14).Value2 = "$/caja 36 atados"
$ws.Cells.Item(429, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(429, 16).Value2 = 361
$ws.Cells.Item(429, 17).Value2 = 36
$ws.Cells.Item(429, 18).Value2 = "Hortaliza"

# Populate the second new row (430) with another new weekly price report.
$ws.Cells.Item(430, 1).Value2  = 4
$ws.Cells.Item(430, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(430, 3).Value2  = "Los Lagos"
$ws.Cells.Item(430, 4).Value2  = 45041
$ws.Cells.Item(430, 5).Value2  = 10
$ws.Cells.Item(430, 6).Value2  = 100112040
$ws.Cells.Item(430, 7).Value2  = "Cilantro"
$ws.Cells.Item(430, 8).Value2  = "Sin especificar"
$ws.Cells.Item(430, 9).Value2  = "Primera"
$ws.Cells.Item(430, 10).Value2 = 160
$ws.Cells.Item(430, 11).Value2 = 5000
$ws.Cells.Item(430, 12).Value2 = 5000
$ws.Cells.Item(430, 13).Value2 = 5000
$ws.Cells.Item(430, 14).Value2 = "$/docena de atados (2 kilos)"
$ws.Cells.Item(430, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(430, 16).Value2 = 2500
$ws.Cells.Item(430, 17).Value2 = 2
$ws.Cells.Item(430, 18).Value2 = "Hortaliza"

# Make sure the D column on the two new rows carries the same date
# number-format style ("s=2") as every other row in this column.
$ws.Range("D429").NumberFormat = $ws.Range("D431").NumberFormat
$ws.Range("D430").NumberFormat = $ws.Range("D431").NumberFormat
